# Generate Report for Handback
#
# This refreshes the handback-status report for the file
# "86bd36d9-8eff-480d-8f23-7f4ce7cbbbb3.md" (row 2 on every sheet):
#   - Overview!G2            -> new "Latest HO Xliff Generate Date"
#   - zh-cn!H2 / zh-cn!K2    -> new Correspond Handoff / Handback datetimes
#   - de-de!H2 / de-de!K2    -> new Correspond Handoff / Handback datetimes
#
# The "ca70f728-..." row (row 3) was already up to date and keeps its
# previous values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 86bd... file
$wsOverview.Range("G2").Value = "2016-08-17 20:49:29"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the 86bd... file
$wsZhCn.Range("H2").Value = "2016-08-17 20:49:24"
$wsZhCn.Range("K2").Value = "2016-08-17 20:49:41"

# de-de sheet: Correspond Handoff / Handback datetimes for the 86bd... file
$wsDeDe.Range("H2").Value = "2016-08-17 20:49:29"
$wsDeDe.Range("K2").Value = "2016-08-17 20:49:48"
